$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find.Execute did not find: $find"
    }
}

# 1) "The monetization of semiconducors started" -> "Industrial semiconducor production started"
Replace-Text "The monetization of semiconducors started" "Industrial semiconducor production started"

# 2) "which did not have the same physical constraints as semiconductor." ->
#    "which does not have the same physical constraints required for manufacturing of semiconductors."
Replace-Text "which did not have the same physical constraints as semiconductor." "which does not have the same physical constraints required for manufacturing of semiconductors."

# 3) "for silicon manufacturing, such as" -> "for silicon work, such as"
Replace-Text "for silicon manufacturing, such as" "for silicon work, such as"

# 4) Append new sentences after "lower Manhattan in New York City."
Replace-Text "lower Manhattan in New York City." "lower Manhattan in New York City. Silicon Valley was initially defined by companies who develop silicon based semiconductors. In recent years the definition has expanded to include software and many other technologies. "

# 5) Replace the whole "These days..." paragraph with the new 2025 paragraph text
Replace-Text "These days people consider the Silicon Valley to include San Francisco, despite the objections of silicon-purists like me. (Silicon and software are not the same)." "This brings us to 2025, where Silicon Valley as a region has expanded to include the entire San Francisco/San Jose region. Some people even apply this term to technology companies no matter where they are located, however many places have their own acroynm containing the word silicon, such as Silicon Alley (NYC), Silcon Prarie (Austin), etc."

# 6) Split the "I enjoyed..." paragraph: insert the new "I have included..." content, then start a
#    fresh paragraph that continues with "The code includes..."
Replace-Text "for a Figure Friday project.  The code includes" "for a Figure Friday project.  I have included the 58 out of 172 companies (34%) in this week’s data set based that are based in the San Francisco bay area, mostly to look at who these companies are, where are they located, and what do they do. With Map Libre, it is easy to see the 3 primary locations where these companies are located: In Santa Clara County (the  original Silicon Valley), in downtown San Francisco, and in San Mateo County. ^pThe code includes"

# 7) Update the screenshot caption text
Replace-Text "showing where most of the companies are located." "showing the locations of Bay Area companies in the dataset."

Write-Output $d.Content.Text
